# Updated symbol list on Mon Jan 23 18:26:36 UTC 2023 with GitHub Actions
# Refresh cryptocurrency Price (column D) and Volume(1h) (column E) values
# for the coin rows on the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the existing (default/unformatted) style so that writing the new
# numeric-looking text values does not change any cell formatting - these
# columns store plain text (e.g. "305.28", "-0.60%"), not real numbers.
$defaultStyle = $ws.Range("D2").Style

$ws.Range("D2").Value = "'305.28"
$ws.Range("D2").Style = $defaultStyle
$ws.Range("E2").Value = "'-0.60%"
$ws.Range("E2").Style = $defaultStyle
$ws.Range("D3").Value = "'36.31"
$ws.Range("D3").Style = $defaultStyle
$ws.Range("E3").Value = "'-1.46%"
$ws.Range("E3").Style = $defaultStyle
$ws.Range("D4").Value = "'5.032"
$ws.Range("D4").Style = $defaultStyle
$ws.Range("E4").Value = "'-0.23%"
$ws.Range("E4").Style = $defaultStyle
$ws.Range("D5").Value = "'0.07923"
$ws.Range("D5").Style = $defaultStyle
$ws.Range("E5").Value = "'0.40%"
$ws.Range("E5").Style = $defaultStyle
$ws.Range("D6").Value = "'2.122"
$ws.Range("D6").Style = $defaultStyle
$ws.Range("E6").Value = "'-3.15%"
$ws.Range("E6").Style = $defaultStyle
$ws.Range("D7").Value = "'7.963"
$ws.Range("D7").Style = $defaultStyle
$ws.Range("E7").Value = "'-0.95%"
$ws.Range("E7").Style = $defaultStyle
$ws.Range("D8").Value = "'4.142"
$ws.Range("D8").Style = $defaultStyle
$ws.Range("E8").Value = "'2.08%"
$ws.Range("E8").Style = $defaultStyle
$ws.Range("D9").Value = "'0.9232"
$ws.Range("D9").Style = $defaultStyle
$ws.Range("E9").Value = "'-0.45%"
$ws.Range("E9").Style = $defaultStyle
$ws.Range("D10").Value = "'0.09737"
$ws.Range("D10").Style = $defaultStyle
$ws.Range("E10").Value = "'-2.04%"
$ws.Range("E10").Style = $defaultStyle
$ws.Range("D11").Value = "'0.1856"
$ws.Range("D11").Style = $defaultStyle
$ws.Range("E11").Value = "'-1.87%"
$ws.Range("E11").Style = $defaultStyle
$ws.Range("D12").Value = "'0.09013"
$ws.Range("D12").Style = $defaultStyle
$ws.Range("E12").Value = "'3.32%"
$ws.Range("E12").Style = $defaultStyle
$ws.Range("D13").Value = "'0.03596"
$ws.Range("D13").Style = $defaultStyle
$ws.Range("E13").Value = "'-0.12%"
$ws.Range("E13").Style = $defaultStyle
$ws.Range("D14").Value = "'0.09915"
$ws.Range("D14").Style = $defaultStyle
$ws.Range("E14").Value = "'-0.50%"
$ws.Range("E14").Style = $defaultStyle
$ws.Range("E15").Value = "'-4.34%"
$ws.Range("E15").Style = $defaultStyle
$ws.Range("D16").Value = "'0.005608"
$ws.Range("D16").Style = $defaultStyle
$ws.Range("E16").Value = "'-0.96%"
$ws.Range("E16").Style = $defaultStyle
$ws.Range("D17").Value = "'3.480"
$ws.Range("D17").Style = $defaultStyle
$ws.Range("E17").Value = "'0.50%"
$ws.Range("E17").Style = $defaultStyle
$ws.Range("D18").Value = "'2.657"
$ws.Range("D18").Style = $defaultStyle
$ws.Range("E18").Value = "'13.70%"
$ws.Range("E18").Style = $defaultStyle
$ws.Range("D19").Value = "'0.3423"
$ws.Range("D19").Style = $defaultStyle
$ws.Range("E19").Value = "'-0.38%"
$ws.Range("E19").Style = $defaultStyle
$ws.Range("D20").Value = "'0.1335"
$ws.Range("D20").Style = $defaultStyle
$ws.Range("E20").Value = "'0.63%"
$ws.Range("E20").Style = $defaultStyle
$ws.Range("D21").Value = "'5.156"
$ws.Range("D21").Style = $defaultStyle
$ws.Range("E21").Value = "'4.58%"
$ws.Range("E21").Style = $defaultStyle
$ws.Range("D22").Value = "'0.2244"
$ws.Range("D22").Style = $defaultStyle
$ws.Range("E22").Value = "'2.01%"
$ws.Range("E22").Style = $defaultStyle
$ws.Range("D23").Value = "'0.04574"
$ws.Range("D23").Style = $defaultStyle
$ws.Range("E23").Value = "'-1.13%"
$ws.Range("E23").Style = $defaultStyle
$ws.Range("D24").Value = "'0.001234"
$ws.Range("D24").Style = $defaultStyle
$ws.Range("E24").Value = "'-1.12%"
$ws.Range("E24").Style = $defaultStyle
$ws.Range("D25").Value = "'0.004816"
$ws.Range("D25").Style = $defaultStyle
$ws.Range("E25").Value = "'-7.92%"
$ws.Range("E25").Style = $defaultStyle
$ws.Range("D26").Value = "'0.0001298"
$ws.Range("D26").Style = $defaultStyle
$ws.Range("E26").Value = "'-7.24%"
$ws.Range("E26").Style = $defaultStyle
$ws.Range("D27").Value = "'0.0004742"
$ws.Range("D27").Style = $defaultStyle
$ws.Range("E27").Value = "'74.55%"
$ws.Range("E27").Style = $defaultStyle
$ws.Range("D39").Value = "'0.01862"
$ws.Range("D39").Style = $defaultStyle
$ws.Range("E39").Value = "'1.34%"
$ws.Range("E39").Style = $defaultStyle
$ws.Range("D40").Value = "'0.04876"
$ws.Range("D40").Style = $defaultStyle
$ws.Range("E40").Value = "'2.05%"
$ws.Range("E40").Style = $defaultStyle
$ws.Range("D41").Value = "'0.007722"
$ws.Range("D41").Style = $defaultStyle
$ws.Range("E41").Value = "'-2.88%"
$ws.Range("E41").Style = $defaultStyle
$ws.Range("D42").Value = "'0.1395"
$ws.Range("D42").Style = $defaultStyle
$ws.Range("E42").Value = "'-1.55%"
$ws.Range("E42").Style = $defaultStyle
$ws.Range("D43").Value = "'0.007726"
$ws.Range("D43").Style = $defaultStyle
$ws.Range("E43").Value = "'1.82%"
$ws.Range("E43").Style = $defaultStyle
$ws.Range("D44").Value = "'0.002297"
$ws.Range("D44").Style = $defaultStyle
$ws.Range("E44").Value = "'4.92%"
$ws.Range("E44").Style = $defaultStyle
$ws.Range("D45").Value = "'0.01121"
$ws.Range("D45").Style = $defaultStyle
$ws.Range("E45").Value = "'10.97%"
$ws.Range("E45").Style = $defaultStyle
$ws.Range("D46").Value = "'0.00006425"
$ws.Range("D46").Style = $defaultStyle
$ws.Range("E46").Value = "'2.61%"
$ws.Range("E46").Style = $defaultStyle
$ws.Range("E47").Value = "'-0.13%"
$ws.Range("E47").Style = $defaultStyle
$ws.Range("E48").Value = "'0.21%"
$ws.Range("E48").Style = $defaultStyle
$ws.Range("D49").Value = "'51.66"
$ws.Range("D49").Style = $defaultStyle
$ws.Range("E49").Value = "'45.06%"
$ws.Range("E49").Style = $defaultStyle
$ws.Range("D50").Value = "'0.001897"
$ws.Range("D50").Style = $defaultStyle
$ws.Range("E50").Value = "'-29.43%"
$ws.Range("E50").Style = $defaultStyle
$ws.Range("D51").Value = "'0.00002097"
$ws.Range("D51").Style = $defaultStyle
$ws.Range("E51").Value = "'-0.13%"
$ws.Range("E51").Style = $defaultStyle

